$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new daily records (rows 445-470, dates 2021-05-03 .. 2021-05-28)
$ws.Cells.Item(445, 1).Value = 44319
$ws.Cells.Item(445, 2).Value = 444
$ws.Cells.Item(445, 3).Value = 837457
$ws.Cells.Item(445, 4).Value = 180
$ws.Cells.Item(445, 5).Value = 23356
$ws.Cells.Item(445, 6).Value = 16977
$ws.Cells.Item(445, 7).Value = 0
$ws.Cells.Item(446, 1).Value = 44320
$ws.Cells.Item(446, 2).Value = 445
$ws.Cells.Item(446, 3).Value = 837715
$ws.Cells.Item(446, 4).Value = 258
$ws.Cells.Item(446, 5).Value = 22833
$ws.Cells.Item(446, 6).Value = 16981
$ws.Cells.Item(446, 7).Value = 4
$ws.Cells.Item(447, 1).Value = 44321
$ws.Cells.Item(447, 2).Value = 446
$ws.Cells.Item(447, 3).Value = 838102
$ws.Cells.Item(447, 4).Value = 387
$ws.Cells.Item(447, 5).Value = 22705
$ws.Cells.Item(447, 6).Value = 16983
$ws.Cells.Item(447, 7).Value = 2
$ws.Cells.Item(448, 1).Value = 44322
$ws.Cells.Item(448, 2).Value = 447
$ws.Cells.Item(448, 3).Value = 838475
$ws.Cells.Item(448, 4).Value = 373
$ws.Cells.Item(448, 5).Value = 22535
$ws.Cells.Item(448, 6).Value = 16988
$ws.Cells.Item(448, 7).Value = 5
$ws.Cells.Item(449, 1).Value = 44323
$ws.Cells.Item(449, 2).Value = 448
$ws.Cells.Item(449, 3).Value = 838852
$ws.Cells.Item(449, 4).Value = 377
$ws.Cells.Item(449, 5).Value = 22421
$ws.Cells.Item(449, 6).Value = 16989
$ws.Cells.Item(449, 7).Value = 1
$ws.Cells.Item(450, 1).Value = 44324
$ws.Cells.Item(450, 2).Value = 449
$ws.Cells.Item(450, 3).Value = 839258
$ws.Cells.Item(450, 4).Value = 406
$ws.Cells.Item(450, 5).Value = 22260
$ws.Cells.Item(450, 6).Value = 16991
$ws.Cells.Item(450, 7).Value = 2
$ws.Cells.Item(451, 1).Value = 44325
$ws.Cells.Item(451, 2).Value = 450
$ws.Cells.Item(451, 3).Value = 839582
$ws.Cells.Item(451, 4).Value = 324
$ws.Cells.Item(451, 5).Value = 22313
$ws.Cells.Item(451, 6).Value = 16992
$ws.Cells.Item(451, 7).Value = 1
$ws.Cells.Item(452, 1).Value = 44326
$ws.Cells.Item(452, 2).Value = 451
$ws.Cells.Item(452, 3).Value = 839740
$ws.Cells.Item(452, 4).Value = 158
$ws.Cells.Item(452, 5).Value = 22102
$ws.Cells.Item(452, 6).Value = 16993
$ws.Cells.Item(452, 7).Value = 1
$ws.Cells.Item(453, 1).Value = 44327
$ws.Cells.Item(453, 2).Value = 452
$ws.Cells.Item(453, 3).Value = 840008
$ws.Cells.Item(453, 4).Value = 268
$ws.Cells.Item(453, 5).Value = 21708
$ws.Cells.Item(453, 6).Value = 16994
$ws.Cells.Item(453, 7).Value = 1
$ws.Cells.Item(454, 1).Value = 44328
$ws.Cells.Item(454, 2).Value = 453
$ws.Cells.Item(454, 3).Value = 840493
$ws.Cells.Item(454, 4).Value = 485
$ws.Cells.Item(454, 5).Value = 21874
$ws.Cells.Item(454, 6).Value = 16998
$ws.Cells.Item(454, 7).Value = 4
$ws.Cells.Item(455, 1).Value = 44329
$ws.Cells.Item(455, 2).Value = 454
$ws.Cells.Item(455, 3).Value = 840929
$ws.Cells.Item(455, 4).Value = 436
$ws.Cells.Item(455, 5).Value = 21969
$ws.Cells.Item(455, 6).Value = 16999
$ws.Cells.Item(455, 7).Value = 1
$ws.Cells.Item(456, 1).Value = 44330
$ws.Cells.Item(456, 2).Value = 455
$ws.Cells.Item(456, 3).Value = 841379
$ws.Cells.Item(456, 4).Value = 450
$ws.Cells.Item(456, 5).Value = 22095
$ws.Cells.Item(456, 6).Value = 16999
$ws.Cells.Item(456, 7).Value = 0
$ws.Cells.Item(457, 1).Value = 44331
$ws.Cells.Item(457, 2).Value = 456
$ws.Cells.Item(457, 3).Value = 841848
$ws.Cells.Item(457, 4).Value = 469
$ws.Cells.Item(457, 5).Value = 22171
$ws.Cells.Item(457, 6).Value = 17006
$ws.Cells.Item(457, 7).Value = 7
$ws.Cells.Item(458, 1).Value = 44332
$ws.Cells.Item(458, 2).Value = 457
$ws.Cells.Item(458, 3).Value = 842182
$ws.Cells.Item(458, 4).Value = 334
$ws.Cells.Item(458, 5).Value = 22275
$ws.Cells.Item(458, 6).Value = 17007
$ws.Cells.Item(458, 7).Value = 1
$ws.Cells.Item(459, 1).Value = 44333
$ws.Cells.Item(459, 2).Value = 458
$ws.Cells.Item(459, 3).Value = 842381
$ws.Cells.Item(459, 4).Value = 199
$ws.Cells.Item(459, 5).Value = 22181
$ws.Cells.Item(459, 6).Value = 17009
$ws.Cells.Item(459, 7).Value = 2
$ws.Cells.Item(460, 1).Value = 44334
$ws.Cells.Item(460, 2).Value = 459
$ws.Cells.Item(460, 3).Value = 842767
$ws.Cells.Item(460, 4).Value = 386
$ws.Cells.Item(460, 5).Value = 21997
$ws.Cells.Item(460, 6).Value = 17011
$ws.Cells.Item(460, 7).Value = 2
$ws.Cells.Item(461, 1).Value = 44335
$ws.Cells.Item(461, 2).Value = 460
$ws.Cells.Item(461, 3).Value = 843278
$ws.Cells.Item(461, 4).Value = 511
$ws.Cells.Item(461, 5).Value = 22089
$ws.Cells.Item(461, 6).Value = 17013
$ws.Cells.Item(461, 7).Value = 2
$ws.Cells.Item(462, 1).Value = 44336
$ws.Cells.Item(462, 2).Value = 461
$ws.Cells.Item(462, 3).Value = 843729
$ws.Cells.Item(462, 4).Value = 451
$ws.Cells.Item(462, 5).Value = 22193
$ws.Cells.Item(462, 6).Value = 17014
$ws.Cells.Item(462, 7).Value = 1
$ws.Cells.Item(463, 1).Value = 44337
$ws.Cells.Item(463, 2).Value = 462
$ws.Cells.Item(463, 3).Value = 844288
$ws.Cells.Item(463, 4).Value = 559
$ws.Cells.Item(463, 5).Value = 22287
$ws.Cells.Item(463, 6).Value = 17017
$ws.Cells.Item(463, 7).Value = 3
$ws.Cells.Item(464, 1).Value = 44338
$ws.Cells.Item(464, 2).Value = 463
$ws.Cells.Item(464, 3).Value = 844811
$ws.Cells.Item(464, 4).Value = 523
$ws.Cells.Item(464, 5).Value = 22328
$ws.Cells.Item(464, 6).Value = 17017
$ws.Cells.Item(464, 7).Value = 0
$ws.Cells.Item(465, 1).Value = 44339
$ws.Cells.Item(465, 2).Value = 464
$ws.Cells.Item(465, 3).Value = 845224
$ws.Cells.Item(465, 4).Value = 413
$ws.Cells.Item(465, 5).Value = 22515
$ws.Cells.Item(465, 6).Value = 17017
$ws.Cells.Item(465, 7).Value = 0
$ws.Cells.Item(466, 1).Value = 44340
$ws.Cells.Item(466, 2).Value = 465
$ws.Cells.Item(466, 3).Value = 845465
$ws.Cells.Item(466, 4).Value = 241
$ws.Cells.Item(466, 5).Value = 22468
$ws.Cells.Item(466, 6).Value = 17018
$ws.Cells.Item(466, 7).Value = 1
$ws.Cells.Item(467, 1).Value = 44341
$ws.Cells.Item(467, 2).Value = 466
$ws.Cells.Item(467, 3).Value = 845840
$ws.Cells.Item(467, 4).Value = 375
$ws.Cells.Item(467, 5).Value = 22171
$ws.Cells.Item(467, 6).Value = 17021
$ws.Cells.Item(467, 7).Value = 3
$ws.Cells.Item(468, 1).Value = 44342
$ws.Cells.Item(468, 2).Value = 467
$ws.Cells.Item(468, 3).Value = 846434
$ws.Cells.Item(468, 4).Value = 594
$ws.Cells.Item(468, 5).Value = 22347
$ws.Cells.Item(468, 6).Value = 17022
$ws.Cells.Item(468, 7).Value = 1
$ws.Cells.Item(469, 1).Value = 44343
$ws.Cells.Item(469, 2).Value = 468
$ws.Cells.Item(469, 3).Value = 847006
$ws.Cells.Item(469, 4).Value = 572
$ws.Cells.Item(469, 5).Value = 22452
$ws.Cells.Item(469, 6).Value = 17022
$ws.Cells.Item(469, 7).Value = 0
$ws.Cells.Item(470, 1).Value = 44344
$ws.Cells.Item(470, 2).Value = 469
$ws.Cells.Item(470, 3).Value = 847604
$ws.Cells.Item(470, 4).Value = 598
$ws.Cells.Item(470, 5).Value = 22534
$ws.Cells.Item(470, 6).Value = 17023
$ws.Cells.Item(470, 7).Value = 1

# Reflect the updated selection/view from the authored workbook
$ws.Range("G2:G470").Select() | Out-Null
